$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 45676
$ws.Range("J109").Value = 45676
$ws.Range("L109").Value = 45676
$ws.Range("N109").Value = -48450
$ws.Range("H116").Value = 7538.3076
$ws.Range("I116").Value = 4350
$ws.Range("J116").Value = 8118
$ws.Range("K116").Value = 4350
$ws.Range("L116").Value = 8118
$ws.Range("M116").Value = -908
$ws.Range("N116").Value = -15002
$ws.Range("H132").Value = 16417.889
$ws.Range("I132").Value = 2409.5925
$ws.Range("J132").Value = 100467.664
$ws.Range("K132").Value = 7228.7775
$ws.Range("L132").Value = 301402.992
$ws.Range("M132").Value = -4698.7775
$ws.Range("N132").Value = -306462.992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 32081.857
$ws.Range("J121").Value = 32081.857
$ws.Range("L121").Value = 32081.857
$ws.Range("N121").Value = -35575.857
$ws.Range("H123").Value = 40878
$ws.Range("J123").Value = 40878
$ws.Range("L123").Value = 40878
$ws.Range("N123").Value = -50678
$ws.Range("H133").Value = 35546.547
$ws.Range("J133").Value = 35546.547
$ws.Range("L133").Value = 35546.547
$ws.Range("N133").Value = -40606.547

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 69958
$ws.Range("J57").Value = 69958
$ws.Range("L57").Value = 69958
$ws.Range("N57").Value = -71398
$ws.Range("H136").Value = 69958
$ws.Range("J136").Value = 69958
$ws.Range("L136").Value = 69958
$ws.Range("N136").Value = -80158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 58000
$ws.Range("J52").Value = 58000
$ws.Range("L52").Value = 58000
$ws.Range("N52").Value = -58588
$ws.Range("H81").Value = 35283
$ws.Range("J81").Value = 35283
$ws.Range("L81").Value = 35283
$ws.Range("N81").Value = -37279
$ws.Range("H82").Value = 36496
$ws.Range("J82").Value = 36496
$ws.Range("L82").Value = 36496
$ws.Range("N82").Value = -37218
$ws.Range("H84").Value = 35283
$ws.Range("J84").Value = 35283
$ws.Range("L84").Value = 105849
$ws.Range("N84").Value = -115833
$ws.Range("H85").Value = 36496
$ws.Range("J85").Value = 36496
$ws.Range("L85").Value = 36496
$ws.Range("N85").Value = -38992
$ws.Range("H88").Value = 35117.5
$ws.Range("J88").Value = 35117.5
$ws.Range("L88").Value = 35117.5
$ws.Range("N88").Value = -35929.5
$ws.Range("H91").Value = 35117.5
$ws.Range("J91").Value = 35117.5
$ws.Range("L91").Value = 35117.5
$ws.Range("N91").Value = -37925.5
$ws.Range("H100").Value = 32844.668
$ws.Range("J100").Value = 39767
$ws.Range("L100").Value = 39767
$ws.Range("N100").Value = -41931
$ws.Range("H111").Value = 25675
$ws.Range("J111").Value = 25675
$ws.Range("L111").Value = 25675
$ws.Range("N111").Value = -33855
$ws.Range("H115").Value = 34256.332
$ws.Range("J115").Value = 34256.332
$ws.Range("L115").Value = 34256.332
$ws.Range("N115").Value = -36606.332
$ws.Range("H119").Value = 34814.668
$ws.Range("J119").Value = 34814.668
$ws.Range("L119").Value = 34814.668
$ws.Range("N119").Value = -44490.668
$ws.Range("H125").Value = 11250
$ws.Range("J125").Value = 11250
$ws.Range("L125").Value = 11250
$ws.Range("N125").Value = -16170
$ws.Range("H131").Value = 35621
$ws.Range("J131").Value = 35621
$ws.Range("L131").Value = 35621
$ws.Range("N131").Value = -45701
$ws.Range("H132").Value = 59139.2
$ws.Range("I132").Value = 1903.3529
$ws.Range("J132").Value = 180765.38
$ws.Range("K132").Value = 5710.0587
$ws.Range("L132").Value = 542296.14
$ws.Range("M132").Value = -3180.0587
$ws.Range("N132").Value = -547356.14
$ws.Range("H134").Value = 585453.9399999999
$ws.Range("I134").Value = 1091.875
$ws.Range("J134").Value = 1754178.1
$ws.Range("K134").Value = 3275.625
$ws.Range("L134").Value = 5262534.300000001
$ws.Range("M134").Value = -740.625
$ws.Range("N134").Value = -5267604.300000001
$ws.Range("H139").Value = 38720
$ws.Range("J139").Value = 37900
$ws.Range("L139").Value = 37900
$ws.Range("N139").Value = -48180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 29661324
$ws.Range("I33").Value = 340
$ws.Range("J33").Value = 44491816
$ws.Range("K33").Value = 2040
$ws.Range("L33").Value = 266950896
$ws.Range("M33").Value = -1757
$ws.Range("N33").Value = -266951462
$ws.Range("H42").Value = 1688.8889
$ws.Range("J42").Value = 650
$ws.Range("L42").Value = 1950
$ws.Range("N42").Value = -3018
$ws.Range("H121").Value = 344604.66
$ws.Range("J121").Value = 775208
$ws.Range("L121").Value = 2325624
$ws.Range("N121").Value = -2328244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 25135.637
$ws.Range("J74").Value = 26149.2
$ws.Range("L74").Value = 26149.2
$ws.Range("N74").Value = -28021.2
$ws.Range("H77").Value = 25135.637
$ws.Range("J77").Value = 26149.2
$ws.Range("L77").Value = 78447.60000000001
$ws.Range("N77").Value = -87807.60000000001
$ws.Range("H100").Value = 36478.6
$ws.Range("J100").Value = 36478.6
$ws.Range("L100").Value = 36478.6
$ws.Range("N100").Value = -38642.6
$ws.Range("H110").Value = 31478.6
$ws.Range("J110").Value = 31478.6
$ws.Range("L110").Value = 31478.6
$ws.Range("N110").Value = -39658.6
$ws.Range("H120").Value = 39313
$ws.Range("J120").Value = 39313
$ws.Range("L120").Value = 39313
$ws.Range("N120").Value = -48989
$ws.Range("H130").Value = 46563.777
$ws.Range("J130").Value = 46563.777
$ws.Range("L130").Value = 46563.777
$ws.Range("N130").Value = -56603.777
$ws.Range("H137").Value = 31520
$ws.Range("J137").Value = 31520
$ws.Range("L137").Value = 31520
$ws.Range("N137").Value = -41720
$ws.Range("H139").Value = 21445.2
$ws.Range("J139").Value = 21445.2
$ws.Range("L139").Value = 21445.2
$ws.Range("N139").Value = -31725.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H88").Value = 20120
$ws.Range("I88").Value = 8085.5
$ws.Range("J88").Value = 44189
$ws.Range("K88").Value = 8085.5
$ws.Range("L88").Value = 44189
$ws.Range("M88").Value = -7657.5
$ws.Range("N88").Value = -45045
$ws.Range("H91").Value = 20120
$ws.Range("I91").Value = 8085.5
$ws.Range("J91").Value = 44189
$ws.Range("K91").Value = 8085.5
$ws.Range("L91").Value = 44189
$ws.Range("M91").Value = -6603.5
$ws.Range("N91").Value = -47153
$ws.Range("H114").Value = 22892.334
$ws.Range("J114").Value = 22892.334
$ws.Range("L114").Value = 22892.334
$ws.Range("N114").Value = -31570.334
$ws.Range("H120").Value = 36930.75
$ws.Range("J120").Value = 36930.75
$ws.Range("L120").Value = 36930.75
$ws.Range("N120").Value = -46606.75
$ws.Range("H121").Value = 20922.666
$ws.Range("J121").Value = 20922.666
$ws.Range("L121").Value = 20922.666
$ws.Range("N121").Value = -24416.666
$ws.Range("H137").Value = 37400
$ws.Range("J137").Value = 37400
$ws.Range("L137").Value = 37400
$ws.Range("N137").Value = -47600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66696.60000000001
$ws.Range("J46").Value = 66696.60000000001
$ws.Range("L46").Value = 66696.60000000001
$ws.Range("N46").Value = -67158.60000000001
$ws.Range("H87").Value = 29900
$ws.Range("J87").Value = 29900
$ws.Range("L87").Value = 29900
$ws.Range("N87").Value = -32396
$ws.Range("H90").Value = 29900
$ws.Range("J90").Value = 29900
$ws.Range("L90").Value = 89700
$ws.Range("N90").Value = -102180
$ws.Range("H93").Value = 36371.43
$ws.Range("J93").Value = 36371.43
$ws.Range("L93").Value = 36371.43
$ws.Range("N93").Value = -41363.43
$ws.Range("H99").Value = 37426.91
$ws.Range("J99").Value = 38675
$ws.Range("L99").Value = 38675
$ws.Range("N99").Value = -44665
$ws.Range("H103").Value = 34416.668
$ws.Range("J103").Value = 34416.668
$ws.Range("L103").Value = 34416.668
$ws.Range("N103").Value = -36760.668
$ws.Range("H106").Value = 32366.666
$ws.Range("J106").Value = 32366.666
$ws.Range("L106").Value = 32366.666
$ws.Range("N106").Value = -34890.666
$ws.Range("H110").Value = 26115.25
$ws.Range("J110").Value = 26115.25
$ws.Range("L110").Value = 26115.25
$ws.Range("N110").Value = -34295.25
$ws.Range("H116").Value = 26052.5
$ws.Range("J116").Value = 26052.5
$ws.Range("L116").Value = 26052.5
$ws.Range("N116").Value = -35230.5
$ws.Range("H128").Value = 46045.668
$ws.Range("J128").Value = 46045.668
$ws.Range("L128").Value = 46045.668
$ws.Range("N128").Value = -56005.668
$ws.Range("H134").Value = 66696.60000000001
$ws.Range("J134").Value = 66696.60000000001
$ws.Range("L134").Value = 200089.8
$ws.Range("N134").Value = -205159.8
